$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.505.19"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +3.34%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.817.43"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +4.42%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.37%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "344.20"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +3.23%  "
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +0.35%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3833"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +2.48%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3547"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +3.81%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "48.98"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -1.46%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.238"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07788"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +3.57%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.001"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +0.33%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "22.43"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +9.39%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.608"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +2.02%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "1.814.92"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "7.197"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +1.61%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001124"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +2.57%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06734"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +0.36%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "86.70"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +2.91%  "
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +0.37%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.65"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +5.09%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.558"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +5.81%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "13.19"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +0.22%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "27.495.93"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +3.49%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.465"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -0.50%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.697"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +7.04%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "22.17"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +12.86%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.469"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +3.86%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "153.99"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +1.04%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.019.37"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +4.60%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "136.33"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +3.19%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.378"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +2.75%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.071"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -1.52%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "13.96"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +6.01%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.08802"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +2.80%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.688"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -2.58%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.637"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +2.60%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.7055"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +12.42%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.06518"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +2.31%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.2260"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +3.32%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.02403"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +1.43%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "9.001"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +3.68%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.300"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +4.35%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "14.88"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +3.24%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.6626"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +8.31%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.9999"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +0.27%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.967"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +1.66%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.195"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +5.45%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "132.52"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +2.74%  "
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -0.07%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "80.94"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +3.83%  "
